$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds two stacked blocks of "energy level" data:
#   Rows 1-13  : A = shared-string 0 ("B"), value in column B   -> alpha orbitals
#   Rows 14-26 : A = shared-string 1 ("R"), value in column C   -> beta orbitals
# plus four label rows further down (HOMO/LUMO captions).
#
# The update trims the last two rows from each 13-row block (the two highest
# energy levels in each series are dropped), which also shifts every row
# below up by two, and then up by two again for the second block -
# collapsing the original blank-row gaps as everything moves closer together.

# Drop the last two rows of the first ("B") block -> rows 12:13 (values
# 1.1285622630755998 / 2.4097166080841657), shifting everything below up by 2.
$ws.Rows("12:13").Delete()

# Drop the last two rows of the second ("R") block, which after the shift
# above now sit at rows 23:24 (values 1.3705312948177686 / 2.4245406670088436).
$ws.Rows("23:24").Delete()

# Close up the blank-row gap before the label rows so B33/B34 (now at
# 31:32) land on 29:30, and C45/C46 (now at 43:44) land on 39:40.
$ws.Rows("31:32").Delete()

# Update the view to match: selection on G22.
$ws.Range("G22").Select()
